$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row was inserted right after the current row 204 (which keeps
# its position but gets a new date). Concretely: row 204's original content
# is duplicated down into a brand-new row 205, pushing every following row
# down by one (old 205 -> new 206, ..., old 317 -> new 318), and row 204's
# Fecha (column D) is updated to the new date.

# 1. Insert a blank row at position 205 (shifts rows 205..317 down to 206..318).
$ws.Rows.Item(205).Insert()

# 2. Duplicate row 204 (its original values) down into the newly inserted row 205.
$src = $ws.Range("A204:T204")
$src.Copy()
$ws.Range("A205").PasteSpecial()
$excel.CutCopyMode = $false

# 3. Update row 204's date (column D) to the new value (Excel serial 44879 = 2022-11-14).
$ws.Range("D204").Value2 = 44879
